$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "1.000", "235.18").
# Force text format on the whole price column first so Excel keeps them as
# literal text (matching the source data) instead of auto-converting to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.182.02"
$ws.Range("E2").Value = "  -1.08%  "

$ws.Range("D3").Value = "1.855.20"
$ws.Range("E3").Value = "  -2.54%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "235.18"
$ws.Range("E5").Value = "  -1.51%  "

$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").Value = "0.4781"
$ws.Range("E7").Value = "  -2.61%  "

$ws.Range("D8").Value = "0.2794"
$ws.Range("E8").Value = "  -4.79%  "

$ws.Range("D9").Value = "0.06423"
$ws.Range("E9").Value = "  -4.25%  "

$ws.Range("D10").Value = "1.857.87"
$ws.Range("E10").Value = "  -2.38%  "

$ws.Range("D11").Value = "0.07374"
$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("D12").Value = "16.13"
$ws.Range("E12").Value = "  -5.34%  "

$ws.Range("D13").Value = "5.089"
$ws.Range("E13").Value = "  -1.92%  "

$ws.Range("D14").Value = "86.75"
$ws.Range("E14").Value = "  -1.65%  "

$ws.Range("D15").Value = "0.6430"
$ws.Range("E15").Value = "  -4.18%  "

$ws.Range("D16").Value = "30.127.57"
$ws.Range("E16").Value = "  -1.16%  "

$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  +0.10%  "

$ws.Range("D18").Value = "13.12"
$ws.Range("E18").Value = "  -2.99%  "

$ws.Range("D19").Value = "0.000007547"
$ws.Range("E19").Value = "  -4.49%  "

$ws.Range("D20").Value = "2.097.21"
$ws.Range("E20").Value = "  -2.33%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "222.40"
$ws.Range("E21").Value = "  +13.92%  "

$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.261"
$ws.Range("E23").Value = "  -4.08%  "

$ws.Range("D24").Value = "6.063"
$ws.Range("E24").Value = "  -1.39%  "

$ws.Range("D25").Value = "9.181"
$ws.Range("E25").Value = "  -3.60%  "

$ws.Range("D26").Value = "163.75"
$ws.Range("E26").Value = "  +0.37%  "

$ws.Range("D27").Value = "18.40"
$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").Value = "1.920"
$ws.Range("E28").Value = "  -1.70%  "

$ws.Range("D29").Value = "1.438"
$ws.Range("E29").Value = "  -2.20%  "

$ws.Range("D30").Value = "0.09159"
$ws.Range("E30").Value = "  -0.38%  "

$ws.Range("D31").Value = "4.224"
$ws.Range("E31").Value = "  -3.03%  "

$ws.Range("D32").Value = "3.930"
$ws.Range("E32").Value = "  -3.91%  "

$ws.Range("D33").Value = "0.04963"
$ws.Range("E33").Value = "  -4.12%  "

$ws.Range("D34").Value = "1.138"
$ws.Range("E34").Value = "  +2.64%  "

$ws.Range("D35").Value = "0.7212"
$ws.Range("E35").Value = "  -3.38%  "

$ws.Range("D36").Value = "2.687"
$ws.Range("E36").Value = "  -1.13%  "

$ws.Range("D37").Value = "0.01827"
$ws.Range("E37").Value = "  +0.47%  "

$ws.Range("E38").Value = "  -3.50%  "

$ws.Range("D39").Value = "0.9000"
$ws.Range("E39").Value = "  -2.82%  "

$ws.Range("D40").Value = "2.029"
$ws.Range("E40").Value = "  -1.72%  "

$ws.Range("D41").Value = "5.884"
$ws.Range("E41").Value = "  -0.78%  "

$ws.Range("D42").Value = "105.37"
$ws.Range("E42").Value = "  -1.44%  "

$ws.Range("D43").Value = "0.9998"
$ws.Range("E43").Value = "  +0.45%  "

$ws.Range("D44").Value = "0.4235"
$ws.Range("E44").Value = "  -3.78%  "

$ws.Range("D45").Value = "0.1303"
$ws.Range("E45").Value = "  -5.30%  "

$ws.Range("D46").Value = "7.251"
$ws.Range("E46").Value = "  -4.87%  "

$ws.Range("D47").Value = "63.86"
$ws.Range("E47").Value = "  -8.12%  "

$ws.Range("D48").Value = "1.495"
$ws.Range("E48").Value = "  +5.88%  "

$ws.Range("D49").Value = "8.654"
$ws.Range("E49").Value = "  -3.91%  "

$ws.Range("D50").Value = "33.55"
$ws.Range("E50").Value = "  -4.35%  "

$ws.Range("D51").Value = "0.05638"
$ws.Range("E51").Value = "  -3.57%  "
